$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 77.127561
$ws.Range("H2").Value = 231.382683
$ws.Range("I2").Value = 0.2899056040435161
$ws.Range("J2").Value = 0.2899056040435161
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 31.642327466673
$ws.Range("R2").Value = 284.780947200057
$ws.Range("S2").Value = 0.001014559413298209
$ws.Range("T2").Value = 0.001014559413298208
$ws.Range("G3").Value = 77.127561
$ws.Range("H3").Value = 231.382683
$ws.Range("I3").Value = 0.2899056040435161
$ws.Range("J3").Value = 0.2899056040435161
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 7859.157553846053
$ws.Range("R3").Value = 70732.41798461448
$ws.Range("S3").Value = 0.2519910169454614
$ws.Range("T3").Value = 0.2519910169454614
$ws.Range("G4").Value = 77.127561
$ws.Range("H4").Value = 231.382683
$ws.Range("I4").Value = 0.2899056040435161
$ws.Range("J4").Value = 0.2899056040435161
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 1150.847100944667
$ws.Range("R4").Value = 10357.623908502
$ws.Range("S4").Value = 0.03690002768475652
$ws.Range("T4").Value = 0.03690002768475652
$ws.Range("I5").Value = 0.443028781054351
$ws.Range("J5").Value = 0.443028781054351
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 48.35526313309389
$ws.Range("R5").Value = 435.197368197845
$ws.Range("S5").Value = 0.001550432326631582
$ws.Range("T5").Value = 0.001550432326631581
$ws.Range("I6").Value = 0.443028781054351
$ws.Range("J6").Value = 0.443028781054351
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.3850883581306573
$ws.Range("T6").Value = 0.3850883581306572
$ws.Range("I7").Value = 0.443028781054351
$ws.Range("J7").Value = 0.443028781054351
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("Q7").Value = 1758.704837713028
$ws.Range("R7").Value = 15828.34353941725
$ws.Range("S7").Value = 0.05638999059706212
$ws.Range("T7").Value = 0.05638999059706211
$ws.Range("I8").Value = 0.267065614902133
$ws.Range("J8").Value = 0.2670656149021329
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 29.14941113229811
$ws.Range("R8").Value = 262.344700190683
$ws.Range("S8").Value = 0.0009346281333925575
$ws.Range("T8").Value = 0.0009346281333925572
$ws.Range("I9").Value = 0.267065614902133
$ws.Range("J9").Value = 0.2670656149021329
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.2321380992698979
$ws.Range("T9").Value = 0.2321380992698978
$ws.Range("I10").Value = 0.267065614902133
$ws.Range("J10").Value = 0.2670656149021329
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("R10").Value = 9541.606507316455
$ws.Range("S10").Value = 0.03399288749884253
$ws.Range("T10").Value = 0.03399288749884252
